# Add a new "2022" column (H) to the table that currently ends at column G
# (years 2018-2021), mirroring the formatting of the existing 2021 column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clone column G's (2021) formatting into column H for the data block
#    (header row 4 through the last data row 37) before writing any values,
#    so every new cell inherits the same number format / font / borders as
#    its neighbour to the left.
$ws.Range("G4:G37").Copy()
$ws.Range("H4:H37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Header cell: 2022
$ws.Range("H4").Value = 2022

# 3. Data values for the 2022 column (rows that are section headers inside
#    the table stay blank, matching column G's existing layout).
$values = @{
    5  = 92.960099223795225
    7  = 96.03949422949897
    8  = 91.012153547624152
    10 = 94.391087218067838
    11 = 91.76755842559642
    13 = 92.942689638142156
    14 = 86.897877953385489
    15 = 96.500794494289821
    16 = 94.135975315309977
    17 = 89.456106196597958
    18 = 94.270923428904894
    19 = 97.027480110114013
    20 = 98.077227596867303
    21 = 90.983384827072243
    23 = 90.468970496790078
    24 = 95.809965597614095
    25 = 88.221110530662017
    27 = 69.811292606515579
    28 = 85.757158930558518
    29 = 93.032103866435918
    30 = 97.325262246493097
    31 = 98.908492141713779
    33 = 91.968006037496949
    34 = 91.809335747904541
    35 = 91.27524653351901
    36 = 94.397288657466234
    37 = 96.740699993405215
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 8).Value = $values[$row]
}

# 4. Match the saved selection: cell I4 (just right of the new column) is
#    what was active when the workbook was last saved.
$ws.Range("I4").Select()
